$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("job_sites")
$wsCity = $wb.Worksheets.Item("city")

# Existing data correction: row 2's id goes from 1 to 0
$ws.Range("A2").Value = 0

# New job-board entry: acumen
$ws.Range("A53").Value = 52
$ws.Range("B53").Value = "acumen"
$ws.Range("C53").Value = "https://www.acumenllc.com/careers.html#accordion-Team1%20.item-1"
$ws.Hyperlinks.Add($ws.Range("C53"), "https://www.acumenllc.com/careers.html", "accordion-Team1%20.item-1")
$ws.Range("D53").Value = "Research & evalation firm for policymaking"
$ws.Range("E53").Value = 20
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 1
$ws.Range("I53").Value = "1,2"

# New city lookup entry: Burlingame (id 14)
$wsCity.Range("A16").Value = 14
$wsCity.Range("B16").Value = "Burlingame"
$wsCity.Range("A17").Select() | Out-Null

# City ids referenced by the new job-board row
$ws.Range("J53").Value = "1,8,14"

# Restore job_sites as the active sheet/selection
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
